# Week 6 Precision and Recall - apply the two text edits described by the
# commit:
#   1. "... a threshold other that 0.5 to make the tradeoff ..."
#        -> "... a threshold other than 0.5 to make the tradeoff ..."
#      (typo fix "that" -> "than"), which also leaves Word's last-edit-
#      position "_GoBack" bookmark sitting right after "than".
#   2. The "_GoBack" bookmark that used to sit inside
#        "... here are precision-recall curves ..."
#      is gone (it moved to edit #1 above), and the two runs that used to
#      be split around it ("precision-recall" + " ") collapse back into a
#      single run "precision-recall ".
#
# NOTE: this runtime coalesces same-formatted sibling runs inside a
# paragraph whenever a Range.Text assignment touches that paragraph, and
# it only stops the coalescing at hard barriers (bookmarks, hyperlinks,
# differing formatting, ...). To reproduce the exact run layout from the
# diff we temporarily drop "fence" bookmarks around the text we edit so
# only the intended runs merge, then remove the fences again.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: " other that 0.5 ..." -> " other than" | bookmark | " 0.5 ..."
# ---------------------------------------------------------------------

$full = $d.Content.Text
$idx = $full.IndexOf("a threshold other that 0.5")

$thresholdEnd = $idx + "a threshold".Length
$thatStart = $thresholdEnd + " other ".Length
$thatEnd = $thatStart + "that".Length

# Fence off "a threshold" from the run we are about to edit so the two
# don't get coalesced into one run.
$d.Bookmarks.Add("_Fence1", $d.Range($thresholdEnd, $thresholdEnd)) | Out-Null

# Fix the typo.
$d.Range($thatStart, $thatEnd).Text = "than"

# Re-seat the (singleton) "_GoBack" bookmark right after "than" - this
# both creates the new bookmark here and removes it from wherever it was
# before (inside the "precision-recall" sentence).
$full = $d.Content.Text
$idx = $full.IndexOf("a threshold other than")
$afterThan = $idx + "a threshold other than".Length
$d.Bookmarks.Add("_GoBack", $d.Range($afterThan, $afterThan)) | Out-Null

# Remove the temporary fence now that the edit is locked in.
$d.Bookmarks("_Fence1").Delete()

# ---------------------------------------------------------------------
# Edit 2: collapse "precision-recall" + " " into one run, now that the
# bookmark that used to separate them is gone.
# ---------------------------------------------------------------------

$full = $d.Content.Text
$idx = $full.IndexOf("For instance, here are precision-recall curves")

$hereAreEnd = $idx + "For instance, here are".Length
$precisionStart = $hereAreEnd + 1
$precisionEnd = $precisionStart + "precision-recall".Length
$spaceEnd = $precisionEnd + 1
$curveEnd = $spaceEnd + "curve".Length
$sEnd = $curveEnd + 1

# Fence off every neighbouring run boundary we want to keep intact so the
# coalescing triggered below only touches "precision-recall" + " ".
$d.Bookmarks.Add("_Fence2", $d.Range($hereAreEnd, $hereAreEnd)) | Out-Null
$d.Bookmarks.Add("_Fence3", $d.Range($precisionStart, $precisionStart)) | Out-Null
$d.Bookmarks.Add("_Fence4", $d.Range($spaceEnd, $spaceEnd)) | Out-Null
$d.Bookmarks.Add("_Fence5", $d.Range($curveEnd, $curveEnd)) | Out-Null
$d.Bookmarks.Add("_Fence6", $d.Range($sEnd, $sEnd)) | Out-Null

# Re-write the space between "precision-recall" and "curve" so the
# runtime merges "precision-recall" with the following space into one run.
# (Round-tripping through a throwaway value forces the coalescing pass -
# assigning the very same text back is treated as a no-op.)
$d.Range($precisionEnd, $spaceEnd).Text = "_"
$d.Range($precisionEnd, $spaceEnd).Text = " "

# Remove the temporary fences - the merge already happened and sticks.
$d.Bookmarks("_Fence2").Delete()
$d.Bookmarks("_Fence3").Delete()
$d.Bookmarks("_Fence4").Delete()
$d.Bookmarks("_Fence5").Delete()
$d.Bookmarks("_Fence6").Delete()

Write-Output "Applied precision/recall edits"
